$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.190551169047839
$ws.Range("D2").Value = 0.2465181130149212

$ws.Range("C3").Value = 0.3615216544701667
$ws.Range("D3").Value = 0.7211584478680613

$ws.Range("C4").Value = 0.1135359368632737
$ws.Range("D4").Value = 0.9106353955660977

$ws.Range("C5").Value = -0.06828372387151962
$ws.Range("D5").Value = 0.9461766046568483

$ws.Range("C6").Value = -1.331382724526328
$ws.Range("D6").Value = 0.1966954009483943

$ws.Range("C7").Value = -1.224488525425866
$ws.Range("D7").Value = 0.2337187523947764

$ws.Range("C8").Value = -0.8202916912524593
$ws.Range("D8").Value = 0.4208491825866674

$ws.Range("C9").Value = -0.2759829897380232
$ws.Range("D9").Value = 0.7851370205798123

$ws.Range("C10").Value = -0.2770837388361497
$ws.Range("D10").Value = 0.7843026387873571

$ws.Range("C11").Value = -0.1314806354323396
$ws.Range("D11").Value = 0.8965901028683476
